$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Component sheet: insert 4 new columns (I:L) before the old "UploadFileLocation"
# column, which slides from I to M, carrying its value/style with it.
$ws2.Range("I1:L1").EntireColumn.Insert()

# --- Row 3: new scenario data (TestCaseName..SubMaterial use existing shared strings,
# the rest are new). Order matters: it reproduces the exact shared-string allocation
# order of the original edit.
$ws2.Range("A3").Value = "Verify that user can add Breakdown for the component"
$ws2.Range("B3").Value = "TC2"
$ws2.Range("C3").Value = "Barbie Barbie Inc."
$ws2.Range("D3").Value = "Automation"
$ws2.Range("E3").Value = "Testing product "
$ws2.Range("F3").Value = "Tons"
$ws2.Range("G3").Value = "Natural Fibre"
$ws2.Range("H3").Value = "Linen Organic"

# --- New header cells for the inserted columns, plus their row-3 data.
$ws2.Range("K1").Value = "Origin"
$ws2.Range("I1").Value = "BreakdownMaterial"
$ws2.Range("I3").Value = "Synthetic Fibre"
$ws2.Range("J1").Value = "BreakdownSubMaterial"
$ws2.Range("J3").Value = "Rubber"
$ws2.Range("L1").Value = "Content"
$ws2.Range("K3").Value = "Algeria"
$ws2.Range("L3").Value = 100

# --- Row 2 loses its one-off formatting (big Menlo font on A2, bold-ish style on C2)
# and the extra row height that went with it.
$ws2.Range("A2").Style = "Normal"
$ws2.Range("C2").Style = "Normal"
$ws2.Rows.Item(2).AutoFit()

# --- Column widths for the newly visible/changed columns.
$ws2.Columns.Item(1).ColumnWidth = 42.166666666666664
$ws2.Columns.Item(3).ColumnWidth = 17.498697916666668
$ws2.Columns.Item(4).ColumnWidth = 24.498697916666668
$ws2.Columns.Item(5).ColumnWidth = 18.666666666666668
$ws2.Columns.Item(7).ColumnWidth = 13.830729166666666
$ws2.Columns.Item(8).ColumnWidth = 14.498697916666666
$ws2.Columns.Item(9).ColumnWidth = 17.998697916666668
$ws2.Columns.Item(10).ColumnWidth = 18.998697916666668

# --- Switch the active sheet/selection: Product moves off D10 onto D3 and loses the
# tab focus; Component becomes the selected tab, scrolled to/selecting J4.
$ws1.Activate()
$ws1.Range("D3").Select()
$ws2.Activate()
$ws2.Range("J4").Select()
